# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 73 (pushing the existing
# rows 73..159 down to 74..160) in the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 73..159 down by one row, creating a blank row 73.
$ws.Rows.Item(73).Insert()

# Populate the newly inserted row 73 with the new record.
$ws.Cells.Item(73, 1).Value  = 4
$ws.Cells.Item(73, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(73, 3).Value  = "Los Lagos"
$ws.Cells.Item(73, 4).Value  = 45174
$ws.Cells.Item(73, 5).Value  = 10
$ws.Cells.Item(73, 6).Value  = 100112026
$ws.Cells.Item(73, 7).Value  = "Haba"
$ws.Cells.Item(73, 8).Value  = "Sin especificar"
$ws.Cells.Item(73, 9).Value  = "Primera"
$ws.Cells.Item(73, 10).Value = 80
$ws.Cells.Item(73, 11).Value = 16000
$ws.Cells.Item(73, 12).Value = 16000
$ws.Cells.Item(73, 13).Value = 16000
$ws.Cells.Item(73, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(73, 15).Value = "Provincia de Limar" + [char]0x00ED
$ws.Cells.Item(73, 16).Value = 640
$ws.Cells.Item(73, 17).Value = 25
$ws.Cells.Item(73, 18).Value = "Hortaliza"
